$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Units")

# --- Data edits: row 2 (Solar_Plant_Kasso) capacity columns ---
# Cap_Output1_existing (J2) gets the 304 that used to (also) live in Cap_Output2_existing (L2)
$ws.Range("J2").Value = 304
$ws.Range("L2").ClearContents()

# --- New minimum_op_point (AH) values ---
$ws.Range("AH2").Value = 0.5
$ws.Range("AH3").Value = 0.5
$ws.Range("AH4").Value = 0.2

# --- Updated "Error messages:" calculated column formula (column AJ) ---
$newFormula = '=IF( Table1[[#This Row],[minimum_op_point]]="", "", IF( COUNTA(Table1[[#This Row],[Cap_Input1_existing]], Table1[[#This Row],[Cap_Input2_existing]], Table1[[#This Row],[Cap_Output1_existing]], Table1[[#This Row],[Cap_Output2_existing]]) = 1, "", IF( COUNTA(Table1[[#This Row],[Cap_Input1_existing]], Table1[[#This Row],[Cap_Input2_existing]], Table1[[#This Row],[Cap_Output1_existing]], Table1[[#This Row],[Cap_Output2_existing]]) = 0, "Capacity missing", "Too many capacities" ) ) )'

for ($r = 2; $r -le 6; $r++) {
    $cell = $ws.Range("AJ$r")
    $cell.Formula = $newFormula
    # Give the formula cell its own (non theme-linked) red font so it gets a
    # dedicated style distinct from the plain header/old-formula style.
    $cell.Font.Name = "Calibri"
}

# --- Row 9 "Please check error message..." helper row is removed ---
$ws.Range("A9").ClearContents()
$ws.Range("B9").ClearContents()

# --- Column AJ is now wider (no longer "best fit") ---
$ws.Columns.Item(36).ColumnWidth = 16.5

# --- Page setup (Page Layout -> Size/Orientation were touched) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection moved ---
[void]$ws.Range("F11").Select()
